# Updated symbol list on Thu Jan 19 18:50:01 UTC 2023 with GitHub Actions
# Refresh of the crypto price/volume snapshot in column D (Price) and
# column E (Volume(1h)) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "291.91"
    "E2" = "-0.78%"
    "E3" = "-0.72%"
    "D4" = "4.893"
    "E4" = "-1.38%"
    "D5" = "0.07273"
    "E5" = "-0.71%"
    "D6" = "2.254"
    "E6" = "25.85%"
    "D7" = "7.678"
    "E7" = "0.23%"
    "D8" = "3.709"
    "E8" = "-1.35%"
    "D9" = "0.8962"
    "E9" = "-1.20%"
    "D10" = "0.1665"
    "E10" = "0.99%"
    "D11" = "0.07923"
    "E11" = "4.79%"
    "D12" = "0.08034"
    "E12" = "-1.02%"
    "E13" = "3.65%"
    "E14" = "0.52%"
    "D15" = "0.001502"
    "E15" = "0.26%"
    "D16" = "0.005860"
    "E16" = "4.44%"
    "D17" = "3.465"
    "E17" = "0.18%"
    "E18" = "-1.01%"
    "D19" = "0.3323"
    "E19" = "1.59%"
    "E20" = "-0.62%"
    "D21" = "3.997"
    "E21" = "-8.42%"
    "D22" = "0.2098"
    "D23" = "0.04526"
    "E23" = "1.21%"
    "E24" = "-1.53%"
    "D25" = "0.004419"
    "E25" = "9.40%"
    "E26" = "3.84%"
    "D27" = "0.0003392"
    "E27" = "-95.48%"
    "D39" = "0.01570"
    "E39" = "-4.51%"
    "D40" = "0.04392"
    "E40" = "-0.01%"
    "D41" = "0.007278"
    "E41" = "-1.16%"
    "D42" = "0.009745"
    "D43" = "0.1316"
    "E43" = "-0.13%"
    "D44" = "0.002028"
    "E44" = "-1.69%"
    "D45" = "0.009360"
    "E45" = "-16.53%"
    "D46" = "0.00005734"
    "E46" = "-4.31%"
    "E47" = "-0.16%"
    "E48" = "5.27%"
    "D49" = "0.002896"
    "E49" = "20.59%"
    "D50" = "0.00002099"
    "E50" = "-0.16%"
    "D51" = "0.0001999"
    "E51" = "-0.16%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force a text number format so Excel keeps the exact literal
    # (leading/trailing zeros, "%" suffix, etc.) as a string instead of
    # silently re-parsing it into a numeric value.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
